$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5
$ws.Cells.Item(5, 4).Value = 45237
$ws.Cells.Item(5, 10).Value = 1000

# Row 6
$ws.Cells.Item(6, 4).Value = 45237
$ws.Cells.Item(6, 9).Value = 'Segunda'
$ws.Cells.Item(6, 11).Value = 1000
$ws.Cells.Item(6, 13).Value = 1000
$ws.Cells.Item(6, 15).Value = 'Región de Ñuble'
$ws.Cells.Item(6, 16).Value = 1000

# Row 7
$ws.Cells.Item(7, 4).Value = 45210
$ws.Cells.Item(7, 10).Value = 200
$ws.Cells.Item(7, 11).Value = 1200
$ws.Cells.Item(7, 12).Value = 1300
$ws.Cells.Item(7, 13).Value = 1250
$ws.Cells.Item(7, 15).Value = 'Región de Ñuble'
$ws.Cells.Item(7, 16).Value = 1250

# Row 8
$ws.Cells.Item(8, 4).Value = 44510
$ws.Cells.Item(8, 10).Value = 600
$ws.Cells.Item(8, 11).Value = 900
$ws.Cells.Item(8, 12).Value = 1000
$ws.Cells.Item(8, 13).Value = 950
$ws.Cells.Item(8, 16).Value = 950

# Row 9
$ws.Cells.Item(9, 4).Value = 44523
$ws.Cells.Item(9, 10).Value = 400
$ws.Cells.Item(9, 11).Value = 800
$ws.Cells.Item(9, 12).Value = 900
$ws.Cells.Item(9, 13).Value = 850
$ws.Cells.Item(9, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(9, 16).Value = 850

# Row 10
$ws.Cells.Item(10, 4).Value = 45203
$ws.Cells.Item(10, 10).Value = 300

# Row 11
$ws.Cells.Item(11, 4).Value = 45203
$ws.Cells.Item(11, 10).Value = 200

# Row 12
$ws.Cells.Item(12, 4).Value = 45205
$ws.Cells.Item(12, 10).Value = 400
$ws.Cells.Item(12, 12).Value = 1300
$ws.Cells.Item(12, 13).Value = 1300
$ws.Cells.Item(12, 16).Value = 1300

# Row 13
$ws.Cells.Item(13, 4).Value = 45205
$ws.Cells.Item(13, 10).Value = 300
$ws.Cells.Item(13, 11).Value = 1500
$ws.Cells.Item(13, 12).Value = 1500
$ws.Cells.Item(13, 13).Value = 1500
$ws.Cells.Item(13, 15).Value = 'Región del Maule'
$ws.Cells.Item(13, 16).Value = 1500

# Row 14
$ws.Cells.Item(14, 4).Value = 44848
$ws.Cells.Item(14, 8).Value = 'Sin especificar'
$ws.Cells.Item(14, 10).Value = 500
$ws.Cells.Item(14, 11).Value = 1300
$ws.Cells.Item(14, 12).Value = 1500
$ws.Cells.Item(14, 13).Value = 1400
$ws.Cells.Item(14, 16).Value = 1400

# Row 15
$ws.Cells.Item(15, 4).Value = 44553
$ws.Cells.Item(15, 10).Value = 8000
$ws.Cells.Item(15, 11).Value = 800
$ws.Cells.Item(15, 12).Value = 900
$ws.Cells.Item(15, 13).Value = 850
$ws.Cells.Item(15, 16).Value = 850

# Row 16
$ws.Cells.Item(16, 4).Value = 44874
$ws.Cells.Item(16, 10).Value = 2000
$ws.Cells.Item(16, 11).Value = 1000
$ws.Cells.Item(16, 12).Value = 1100
$ws.Cells.Item(16, 13).Value = 1050
$ws.Cells.Item(16, 16).Value = 1050

# Row 17
$ws.Cells.Item(17, 4).Value = 44504
$ws.Cells.Item(17, 10).Value = 500
$ws.Cells.Item(17, 13).Value = 950
$ws.Cells.Item(17, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(17, 16).Value = 950

# Row 18
$ws.Cells.Item(18, 4).Value = 44882
$ws.Cells.Item(18, 10).Value = 1200
$ws.Cells.Item(18, 11).Value = 1000
$ws.Cells.Item(18, 12).Value = 1100
$ws.Cells.Item(18, 13).Value = 1050
$ws.Cells.Item(18, 15).Value = 'Región de Ñuble'
$ws.Cells.Item(18, 16).Value = 1050

# Row 19
$ws.Cells.Item(19, 4).Value = 45202
$ws.Cells.Item(19, 10).Value = 300
$ws.Cells.Item(19, 15).Value = 'Provincia de Diguillín'

# Row 20
$ws.Cells.Item(20, 9).Value = 'Primera'
$ws.Cells.Item(20, 11).Value = 1500
$ws.Cells.Item(20, 12).Value = 1500
$ws.Cells.Item(20, 13).Value = 1500
$ws.Cells.Item(20, 16).Value = 1500

# Row 21
$ws.Cells.Item(21, 4).Value = 45191
$ws.Cells.Item(21, 9).Value = 'Segunda'
$ws.Cells.Item(21, 10).Value = 200
$ws.Cells.Item(21, 11).Value = 1300
$ws.Cells.Item(21, 12).Value = 1300
$ws.Cells.Item(21, 13).Value = 1300
$ws.Cells.Item(21, 15).Value = 'Región del Maule'
$ws.Cells.Item(21, 16).Value = 1300

# Row 22
$ws.Cells.Item(22, 4).Value = 44532
$ws.Cells.Item(22, 10).Value = 240
$ws.Cells.Item(22, 11).Value = 800
$ws.Cells.Item(22, 12).Value = 900
$ws.Cells.Item(22, 13).Value = 850
$ws.Cells.Item(22, 16).Value = 850

# Row 23
$ws.Cells.Item(23, 4).Value = 45215
$ws.Cells.Item(23, 10).Value = 450
$ws.Cells.Item(23, 11).Value = 1300
$ws.Cells.Item(23, 12).Value = 1500
$ws.Cells.Item(23, 13).Value = 1389
$ws.Cells.Item(23, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(23, 16).Value = 1389

# Row 24
$ws.Cells.Item(24, 4).Value = 44900

# Row 25
$ws.Cells.Item(25, 4).Value = 44894
$ws.Cells.Item(25, 10).Value = 1200
$ws.Cells.Item(25, 11).Value = 900
$ws.Cells.Item(25, 12).Value = 1000
$ws.Cells.Item(25, 13).Value = 950
$ws.Cells.Item(25, 15).Value = 'Región de Ñuble'
$ws.Cells.Item(25, 16).Value = 950

# Row 26
$ws.Cells.Item(26, 4).Value = 44517
$ws.Cells.Item(26, 10).Value = 500
$ws.Cells.Item(26, 13).Value = 850
$ws.Cells.Item(26, 16).Value = 850

# Row 27
$ws.Cells.Item(27, 4).Value = 44505
$ws.Cells.Item(27, 10).Value = 440
$ws.Cells.Item(27, 11).Value = 900
$ws.Cells.Item(27, 12).Value = 1000
$ws.Cells.Item(27, 13).Value = 950
$ws.Cells.Item(27, 16).Value = 950

# Row 28
$ws.Cells.Item(28, 4).Value = 44524

# Row 29
$ws.Cells.Item(29, 4).Value = 44518
$ws.Cells.Item(29, 10).Value = 400

# Row 30
$ws.Cells.Item(30, 4).Value = 44545
$ws.Cells.Item(30, 10).Value = 4000
$ws.Cells.Item(30, 11).Value = 800
$ws.Cells.Item(30, 12).Value = 900
$ws.Cells.Item(30, 13).Value = 850
$ws.Cells.Item(30, 16).Value = 850

# Row 31
$ws.Cells.Item(31, 4).Value = 44831
$ws.Cells.Item(31, 10).Value = 200
$ws.Cells.Item(31, 11).Value = 2000
$ws.Cells.Item(31, 12).Value = 2200
$ws.Cells.Item(31, 13).Value = 2100
$ws.Cells.Item(31, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(31, 16).Value = 2100

# Row 32
$ws.Cells.Item(32, 4).Value = 45225
$ws.Cells.Item(32, 10).Value = 500
$ws.Cells.Item(32, 11).Value = 1200
$ws.Cells.Item(32, 12).Value = 1200
$ws.Cells.Item(32, 13).Value = 1200
$ws.Cells.Item(32, 15).Value = 'Región de Ñuble'
$ws.Cells.Item(32, 16).Value = 1200

# Row 33
$ws.Cells.Item(33, 4).Value = 45218
$ws.Cells.Item(33, 10).Value = 300
$ws.Cells.Item(33, 11).Value = 1300
$ws.Cells.Item(33, 12).Value = 1300
$ws.Cells.Item(33, 13).Value = 1300
$ws.Cells.Item(33, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(33, 16).Value = 1300

# Row 34
$ws.Cells.Item(34, 4).Value = 44895
$ws.Cells.Item(34, 10).Value = 1200
$ws.Cells.Item(34, 11).Value = 900
$ws.Cells.Item(34, 12).Value = 1000
$ws.Cells.Item(34, 13).Value = 950
$ws.Cells.Item(34, 15).Value = 'Región de Ñuble'
$ws.Cells.Item(34, 16).Value = 950

# Row 35
$ws.Cells.Item(35, 4).Value = 45212
$ws.Cells.Item(35, 10).Value = 600
$ws.Cells.Item(35, 11).Value = 1200
$ws.Cells.Item(35, 12).Value = 1300
$ws.Cells.Item(35, 13).Value = 1250
$ws.Cells.Item(35, 16).Value = 1250

# Row 36
$ws.Cells.Item(36, 4).Value = 44890
$ws.Cells.Item(36, 10).Value = 160

# Row 37
$ws.Cells.Item(37, 4).Value = 44503
$ws.Cells.Item(37, 10).Value = 400

# Row 38
$ws.Cells.Item(38, 4).Value = 44910
$ws.Cells.Item(38, 10).Value = 1200

# Row 39
$ws.Cells.Item(39, 4).Value = 44516
$ws.Cells.Item(39, 10).Value = 400
$ws.Cells.Item(39, 11).Value = 900
$ws.Cells.Item(39, 12).Value = 1000
$ws.Cells.Item(39, 13).Value = 950
$ws.Cells.Item(39, 16).Value = 950
